# Quality Improvement Plans.docx edit
#
# This script rewrites the second paragraph (the bullet list that starts
# with "- Tighten up input validation...") so that it:
#   1. Splits the first bullet's run around "GiGo" and wraps it with
#      w:proofErr spellStart/spellEnd markers (as Word's spell checker
#      would when it flags "GiGo" as a misspelling).
#   2. Adds a new bullet "- Use custom validators to improve white space
#      input issues" right after the first bullet.
#   3. Moves the "_GoBack" bookmark so it now sits right after the new
#      "white space input issues" bullet (instead of after the old
#      "Phase 3 pagination suggestion" bullet).
#   4. Adds a new trailing bullet "- Use pipes to better format large
#      numbers " after the "Phase 3 pagination suggestion" bullet.
#
# Because the target paragraph's content (runs, proofErr marks and the
# bookmark) all change together, the whole paragraph body is rebuilt in
# one shot via Range.InsertXML - this avoids the interop layer's flaky
# handling of InsertXML on ranges obtained from a live Find object, and
# keeps every run's formatting explicit.

$d = $word.ActiveDocument

# --- locate the target paragraph (bullet list) without hard-coding its index ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("- Tighten up input validation")) {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate the 'Tighten up input validation' paragraph"
}

$startPos = $target.Range.Start
$endPos = $target.Range.End

# Re-wrap the same bounds in a plain Range (not reused from Find) before
# calling InsertXML, since this interop's InsertXML mis-targets the
# insertion point when called on a range that was just populated by
# Find.Execute.
$targetRange = $d.Range($startPos, $endPos)

# --- build the replacement paragraph content --------------------------------
$lang = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'

$inner =
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  "<w:r>$lang<w:t xml:space=`"preserve`">- Tighten up input validation to reduce </w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  "<w:r>$lang<w:t>GiGo</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r>$lang<w:t>, increase security and (hopefully) eliminate program failure due to unforeseen inputs.</w:t></w:r>" +
  "<w:r>$lang<w:t xml:space=`"preserve`"> </w:t></w:r>" +
  "<w:r>$lang<w:br/><w:t>- Use custom validators to improve white space input issues</w:t></w:r>" +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  "<w:r>$lang<w:br/><w:t>- Improve response to user actions.</w:t></w:r>" +
  "<w:r>$lang<w:br/><w:t>- Clean up and improve UI aesthetics.</w:t></w:r>" +
  "<w:r w:rsidR=`"00420348`">$lang<w:br/><w:t>- Phase 3 pagination suggestion</w:t></w:r>" +
  "<w:r>$lang<w:br/><w:t xml:space=`"preserve`">- Use pipes to better format large numbers </w:t></w:r>" +
  "<w:r>$lang<w:br/></w:r>"

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' + $inner + '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$targetRange.InsertXML($xmlFrag)
